# Update gh-pages generated output (456a3b4)
# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3252
$wsExpo.Range("F3").Value = 3
$wsExpo.Range("F4").Value = 54
$wsExpo.Range("F5").Value = 1161
$wsExpo.Range("F6").Value = 308

# Sheet "全部类型" (all types) mirrors the same data plus one extra row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3252
$wsAll.Range("F3").Value = 3
$wsAll.Range("F4").Value = 54
$wsAll.Range("F5").Value = 1161
$wsAll.Range("F7").Value = 308
